# Use "Then" instead of "Assert" to match BDD syntax (Given/When/Then).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell that previously held the literal "Assert" keyword.
$ws.Range("A12").Value = "Then"

# The sheet has a conditional-formatting rule that highlights cells equal to
# the literal "Assert" keyword; update it to match the new "Then" keyword.
$cfRange = $ws.Range("A1:XFD1048576")
for ($i = 1; $i -le $cfRange.FormatConditions.Count; $i++) {
    $fc = $cfRange.FormatConditions.Item($i)
    if ($fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
        break
    }
}
